$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, '三花智控', '三花智控', '远大控股'),
    @(3, '大有能源', '士兰微', '大有能源'),
    @(4, '士兰微', '闻泰科技', '山子高科'),
    @(5, '华天科技', '山子高科', '三花智控'),
    @(6, '山子高科', '大有能源', '合肥城建'),
    @(7, '卧龙电驱', '东方财富', '科大讯飞'),
    @(8, '白银有色', '华天科技', '宝泰隆'),
    @(9, '宝泰隆', '卧龙电驱', '首开股份'),
    @(10, '远大控股', '东信和平', '卧龙电驱'),
    @(11, '合肥城建', '宝泰隆', '楚江新材'),
    @(12, '黄河旋风', '海峡股份', '万润科技'),
    @(13, '郑州煤电', '黄河旋风', '华天科技'),
    @(14, '首开股份', '大众公用', '白银有色'),
    @(15, '中际旭创', '郑州煤电', '黄河旋风'),
    @(16, '创新医疗', '安泰集团', '士兰微'),
    @(17, '闻泰科技', '合肥城建', '工业富联'),
    @(18, '大众公用', '四方达', '郑州煤电'),
    @(19, '东信和平', '远大控股', '国新能源'),
    @(20, '万润科技', '首开股份', '海通发展'),
    @(21, '海峡股份', '常山北明', '青山纸业')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

$wb.Save()
